# Correção nos dados e inicio da analise PNAD 2009
# The "grandes regiões e unidades da federação" header row (row 6, which had
# no data of its own) is removed entirely. Excel shifts every row below it
# up by one, and the now-unused shared string gets dropped automatically
# when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
